$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three obsolete "ECs" sending-cluster rows (old rows 2-4).
# This shifts old rows 5-10 up to become rows 2-7, preserving columns A-D,
# and leaves only 6 data rows + 1 header row (dimension A1:T7).
$ws.Rows(2).EntireRow.Delete() | Out-Null
$ws.Rows(2).EntireRow.Delete() | Out-Null
$ws.Rows(2).EntireRow.Delete() | Out-Null

# Refresh the remaining rows with the new TPM-derived values.
# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Col2a1"
$ws.Cells.Item(2,3).Value = "Ddr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.2021233333333333
$ws.Cells.Item(2,8).Value = 0.60637
$ws.Cells.Item(2,9).Value = 0.9764223557676824
$ws.Cells.Item(2,10).Value = 0.9764223557676823
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.366183
$ws.Cells.Item(2,14).Value = 1.098549
$ws.Cells.Item(2,15).Value = 0.0639836884691917
$ws.Cells.Item(2,16).Value = 0.0639836884691917
$ws.Cells.Item(2,17).Value = 0.07401412856999999
$ws.Cells.Item(2,18).Value = 0.6661271571299999
$ws.Cells.Item(2,19).Value = 0.06247510382579365
$ws.Cells.Item(2,20).Value = 0.06247510382579365

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Col2a1"
$ws.Cells.Item(3,3).Value = "Ddr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.2021233333333333
$ws.Cells.Item(3,8).Value = 0.60637
$ws.Cells.Item(3,9).Value = 0.9764223557676824
$ws.Cells.Item(3,10).Value = 0.9764223557676823
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.589504333333333
$ws.Cells.Item(3,14).Value = 4.768513
$ws.Cells.Item(3,15).Value = 0.2777364052521014
$ws.Cells.Item(3,16).Value = 0.2777364052521014
$ws.Cells.Item(3,17).Value = 0.3212759142011111
$ws.Cells.Item(3,18).Value = 2.89148322781
$ws.Cells.Item(3,19).Value = 0.2711880350987046
$ws.Cells.Item(3,20).Value = 0.2711880350987045

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Col2a1"
$ws.Cells.Item(4,3).Value = "Ddr1"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.2021233333333333
$ws.Cells.Item(4,8).Value = 0.60637
$ws.Cells.Item(4,9).Value = 0.9764223557676824
$ws.Cells.Item(4,10).Value = 0.9764223557676823
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.767380666666666
$ws.Cells.Item(4,14).Value = 11.302142
$ws.Cells.Item(4,15).Value = 0.6582799062787069
$ws.Cells.Item(4,16).Value = 0.6582799062787069
$ws.Cells.Item(4,17).Value = 0.7614755382822221
$ws.Cells.Item(4,18).Value = 6.853279844539999
$ws.Cells.Item(4,19).Value = 0.6427592168431843
$ws.Cells.Item(4,20).Value = 0.6427592168431842

# Row 5
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,2).Value = "Col2a1"
$ws.Cells.Item(5,3).Value = "Ddr1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.004880666666666667
$ws.Cells.Item(5,8).Value = 0.014642
$ws.Cells.Item(5,9).Value = 0.02357764423231757
$ws.Cells.Item(5,10).Value = 0.02357764423231757
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.366183
$ws.Cells.Item(5,14).Value = 1.098549
$ws.Cells.Item(5,15).Value = 0.0639836884691917
$ws.Cells.Item(5,16).Value = 0.0639836884691917
$ws.Cells.Item(5,17).Value = 0.001787217162
$ws.Cells.Item(5,18).Value = 0.016084954458
$ws.Cells.Item(5,19).Value = 0.001508584643398042
$ws.Cells.Item(5,20).Value = 0.001508584643398042

# Row 6
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Col2a1"
$ws.Cells.Item(6,3).Value = "Ddr1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.004880666666666667
$ws.Cells.Item(6,8).Value = 0.014642
$ws.Cells.Item(6,9).Value = 0.02357764423231757
$ws.Cells.Item(6,10).Value = 0.02357764423231757
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.589504333333333
$ws.Cells.Item(6,14).Value = 4.768513
$ws.Cells.Item(6,15).Value = 0.2777364052521014
$ws.Cells.Item(6,16).Value = 0.2777364052521014
$ws.Cells.Item(6,17).Value = 0.007757840816222223
$ws.Cells.Item(6,18).Value = 0.069820567346
$ws.Cells.Item(6,19).Value = 0.006548370153396825
$ws.Cells.Item(6,20).Value = 0.006548370153396825

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Col2a1"
$ws.Cells.Item(7,3).Value = "Ddr1"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.004880666666666667
$ws.Cells.Item(7,8).Value = 0.014642
$ws.Cells.Item(7,9).Value = 0.02357764423231757
$ws.Cells.Item(7,10).Value = 0.02357764423231757
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.767380666666666
$ws.Cells.Item(7,14).Value = 11.302142
$ws.Cells.Item(7,15).Value = 0.6582799062787069
$ws.Cells.Item(7,16).Value = 0.6582799062787069
$ws.Cells.Item(7,17).Value = 0.01838732924044444
$ws.Cells.Item(7,18).Value = 0.165485963164
$ws.Cells.Item(7,19).Value = 0.01552068943552271
$ws.Cells.Item(7,20).Value = 0.01552068943552271

